$d = $word.ActiveDocument

$d.Content.Find.Execute("769×5=", $true, $false, $false, $false, $false, $true, 1, $false, "274×7=", 2) | Out-Null
$d.Content.Find.Execute("691×9=", $true, $false, $false, $false, $false, $true, 1, $false, "145×2=", 2) | Out-Null
$d.Content.Find.Execute("351×9=", $true, $false, $false, $false, $false, $true, 1, $false, "654×5=", 2) | Out-Null
$d.Content.Find.Execute("241×7=", $true, $false, $false, $false, $false, $true, 1, $false, "555×4=", 2) | Out-Null
$d.Content.Find.Execute("356×7=", $true, $false, $false, $false, $false, $true, 1, $false, "531×7=", 2) | Out-Null
$d.Content.Find.Execute("528×2=", $true, $false, $false, $false, $false, $true, 1, $false, "621×7=", 2) | Out-Null
$d.Content.Find.Execute("269×5=", $true, $false, $false, $false, $false, $true, 1, $false, "791×4=", 2) | Out-Null
$d.Content.Find.Execute("275×7=", $true, $false, $false, $false, $false, $true, 1, $false, "101×2=", 2) | Out-Null
$d.Content.Find.Execute("778×3=", $true, $false, $false, $false, $false, $true, 1, $false, "415×5=", 2) | Out-Null
$d.Content.Find.Execute("117×2=", $true, $false, $false, $false, $false, $true, 1, $false, "640×6=", 2) | Out-Null
$d.Content.Find.Execute("356×8=", $true, $false, $false, $false, $false, $true, 1, $false, "302×2=", 2) | Out-Null
$d.Content.Find.Execute("872×3=", $true, $false, $false, $false, $false, $true, 1, $false, "493×4=", 2) | Out-Null
$d.Content.Find.Execute("768×3=", $true, $false, $false, $false, $false, $true, 1, $false, "447×9=", 2) | Out-Null
$d.Content.Find.Execute("491×9=", $true, $false, $false, $false, $false, $true, 1, $false, "112×7=", 2) | Out-Null
$d.Content.Find.Execute("541×7=", $true, $false, $false, $false, $false, $true, 1, $false, "627×8=", 2) | Out-Null
$d.Content.Find.Execute("902×3=", $true, $false, $false, $false, $false, $true, 1, $false, "902×6=", 2) | Out-Null
$d.Content.Find.Execute("599×8=", $true, $false, $false, $false, $false, $true, 1, $false, "367×8=", 2) | Out-Null
$d.Content.Find.Execute("834×3=", $true, $false, $false, $false, $false, $true, 1, $false, "190×7=", 2) | Out-Null
$d.Content.Find.Execute("326×6=", $true, $false, $false, $false, $false, $true, 1, $false, "432×6=", 2) | Out-Null
$d.Content.Find.Execute("672×7=", $true, $false, $false, $false, $false, $true, 1, $false, "143×8=", 2) | Out-Null
$d.Content.Find.Execute("456×6=", $true, $false, $false, $false, $false, $true, 1, $false, "565×4=", 2) | Out-Null
$d.Content.Find.Execute("668×6=", $true, $false, $false, $false, $false, $true, 1, $false, "476×4=", 2) | Out-Null
$d.Content.Find.Execute("813×4=", $true, $false, $false, $false, $false, $true, 1, $false, "556×4=", 2) | Out-Null
$d.Content.Find.Execute("445×7=", $true, $false, $false, $false, $false, $true, 1, $false, "231×9=", 2) | Out-Null
$d.Content.Find.Execute("740×4=", $true, $false, $false, $false, $false, $true, 1, $false, "137×3=", 2) | Out-Null
